$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 329; this shifts the existing rows 329-370
# down to 330-371 (matching the diff, which re-numbers all subsequent rows).
$ws.Rows("329:329").Insert()

# --- Fill in the new row 329 with the new dataset entry ---
# Order matters for the shared-strings table: new distinct strings are
# appended in first-use order, so we set A, then C (hyperlink), then D,
# to reproduce the exact new shared-string index assignment from the diff
# (1517 = name, 1518 = link URL, 1519 = topics).

# A329: name
$ws.Range("A329").Value = "The Comparative Legislators Database"

# B329: category (reuses existing shared string)
$ws.Range("B329").Value = "parties and politicians"

# C329: link -- add as a real hyperlink (this also sets the cell text to
# the URL since the cell is currently empty), then restore the "Link"
# cell style so it matches the style used by the other link cells.
$linkCell = $ws.Range("C329")
$null = $ws.Hyperlinks.Add($linkCell, "https://github.com/saschagobel/legislatoR")
$linkCell.Style = "Link"

# D329: topics
$ws.Range("D329").Value = "legislators, politicians"

# F329:J329: region flags (africa, asia, easteurope, latinamerica, westeurope)
$ws.Range("F329").Value = 0
$ws.Range("G329").Value = 0
$ws.Range("H329").Value = 1
$ws.Range("I329").Value = 0
$ws.Range("J329").Value = 1

# K329/L329: year_start / year_end
$ws.Range("K329").Value = 1789
$ws.Range("L329").Value = 2019

# M329/N329: availability / registration (reuse existing shared strings)
$ws.Range("M329").Value = "online"
$ws.Range("N329").Value = "no"

# O329: free
$ws.Range("O329").Value = 1

# AB329: revised (plain number, yyyymmdd)
$ws.Range("AB329").Value = 20191129

# Restore the sheet's selection/view to where the author ended up after
# editing (no custom topLeftCell, cursor parked on A13).
$ws.Range("A13").Select()
